{"js": "// \"Versi\" + \"on\" -> merge into a single \"Version\" run (text unchanged,\n// but re-typing it collapses the two adjacent runs Word had split).\nconst body = context.document.body;\n\nlet results = body.search(\"Version\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"Version\", Word.InsertLocation.replace);\nawait context.sync();\n\n// \" 2\" + \".\" -> \" 1.\" (spans the old \" 2\" run and the trailing \".\" run,\n// collapsing them into a single run and bumping the version number).\nresults = body.search(\" 2.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Versi\" + \"on\" are two adjacent runs that spell \"Version\". Re-typing\n#    the word over itself collapses them into a single run (no visible\n#    text change).\n$rVersion = $d.Range()\n$fVersion = $rVersion.Find\n$fVersion.Text = \"Version\"\n$fVersion.Replacement.Text = \"Version\"\n$fVersion.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 2) Drop the final \".\" -- it lives in its own trailing run, entirely\n#    AFTER the (empty) \"_GoBack\" bookmark, so deleting it doesn't touch\n#    the bookmark's position.\n$rPeriod = $d.Range()\n$fPeriod = $rPeriod.Find\n$fPeriod.Text = \".\"\n$fPeriod.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null) | Out-Null\n$rPeriod.Text = \"\"\n\n# 3) Bump the version number: \"2\" -> \"1.\" -- this run sits entirely\n#    BEFORE the bookmark, so this also leaves the bookmark undisturbed.\n$rNum = $d.Range()\n$fNum = $rNum.Find\n$fNum.Text = \"2\"\n$fNum.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null) | Out-Null\n$rNum.Text = \"1.\"\n"}
